{"js": "// Update the three-digit x one-digit multiplication prompts in the\n// practice-sheet table. Every \"<old>=\" expression in the document is\n// replaced by its corresponding \"<new>=\" expression (1:1 unique mapping).\nconst replacements = [\n  [\"266\u00d77=\", \"642\u00d73=\"],\n  [\"423\u00d75=\", \"643\u00d74=\"],\n  [\"927\u00d74=\", \"683\u00d78=\"],\n  [\"729\u00d78=\", \"392\u00d75=\"],\n  [\"788\u00d79=\", \"917\u00d73=\"],\n  [\"561\u00d76=\", \"773\u00d78=\"],\n  [\"820\u00d75=\", \"298\u00d78=\"],\n  [\"105\u00d78=\", \"137\u00d76=\"],\n  [\"105\u00d72=\", \"354\u00d72=\"],\n  [\"386\u00d75=\", \"565\u00d72=\"],\n  [\"933\u00d74=\", \"841\u00d74=\"],\n  [\"297\u00d73=\", \"347\u00d76=\"],\n  [\"870\u00d72=\", \"899\u00d78=\"],\n  [\"947\u00d72=\", \"298\u00d74=\"],\n  [\"341\u00d77=\", \"369\u00d75=\"],\n  [\"444\u00d79=\", \"393\u00d78=\"],\n  [\"629\u00d73=\", \"792\u00d79=\"],\n  [\"755\u00d78=\", \"434\u00d74=\"],\n  [\"688\u00d73=\", \"481\u00d78=\"],\n  [\"583\u00d72=\", \"763\u00d76=\"],\n  [\"863\u00d72=\", \"508\u00d74=\"],\n  [\"713\u00d75=\", \"436\u00d78=\"],\n  [\"192\u00d77=\", \"142\u00d75=\"],\n  [\"275\u00d73=\", \"590\u00d77=\"],\n  [\"531\u00d79=\", \"290\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit x one-digit multiplication prompts in the\n# practice-sheet table. Every \"<old>=\" expression in the document is\n# replaced by its corresponding \"<new>=\" expression (1:1 unique mapping).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"266\u00d77=\", \"642\u00d73=\"),\n  @(\"423\u00d75=\", \"643\u00d74=\"),\n  @(\"927\u00d74=\", \"683\u00d78=\"),\n  @(\"729\u00d78=\", \"392\u00d75=\"),\n  @(\"788\u00d79=\", \"917\u00d73=\"),\n  @(\"561\u00d76=\", \"773\u00d78=\"),\n  @(\"820\u00d75=\", \"298\u00d78=\"),\n  @(\"105\u00d78=\", \"137\u00d76=\"),\n  @(\"105\u00d72=\", \"354\u00d72=\"),\n  @(\"386\u00d75=\", \"565\u00d72=\"),\n  @(\"933\u00d74=\", \"841\u00d74=\"),\n  @(\"297\u00d73=\", \"347\u00d76=\"),\n  @(\"870\u00d72=\", \"899\u00d78=\"),\n  @(\"947\u00d72=\", \"298\u00d74=\"),\n  @(\"341\u00d77=\", \"369\u00d75=\"),\n  @(\"444\u00d79=\", \"393\u00d78=\"),\n  @(\"629\u00d73=\", \"792\u00d79=\"),\n  @(\"755\u00d78=\", \"434\u00d74=\"),\n  @(\"688\u00d73=\", \"481\u00d78=\"),\n  @(\"583\u00d72=\", \"763\u00d76=\"),\n  @(\"863\u00d72=\", \"508\u00d74=\"),\n  @(\"713\u00d75=\", \"436\u00d78=\"),\n  @(\"192\u00d77=\", \"142\u00d75=\"),\n  @(\"275\u00d73=\", \"590\u00d77=\"),\n  @(\"531\u00d79=\", \"290\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
